$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("B2").Value = 1.02
$ws.Range("C2").Value = 1.02147381031754
$ws.Range("D2").Value = 1.026935930692248
$ws.Range("E2").Value = 0.9926147277508489
$ws.Range("F2").Value = 1.032195857781933
$ws.Range("I2").Value = 1.029693403081173
$ws.Range("J2").Value = 1.026665076591334
$ws.Range("K2").Value = 1.029757124067765
$ws.Range("L2").Value = 0.9955398523336033
$ws.Range("M2").Value = 1.035001787395342
$ws.Range("N2").Value = 1.012885763548367
# Row 3
$ws.Range("B3").Value = 1.02
$ws.Range("C3").Value = 1.02267872641413
$ws.Range("D3").Value = 1.02783275713465
$ws.Range("E3").Value = 0.9936372048519304
$ws.Range("F3").Value = 1.033565685585866
$ws.Range("I3").Value = 1.029934320392649
$ws.Range("J3").Value = 1.027506297678012
$ws.Range("K3").Value = 1.030461328387339
$ws.Range("L3").Value = 0.9963617723202692
$ws.Range("M3").Value = 1.036178835819625
$ws.Range("N3").Value = 1.013170150947261
# Row 4
$ws.Range("B4").Value = 1.02
$ws.Range("C4").Value = 1.02345761553967
$ws.Range("D4").Value = 1.028412227140846
$ws.Range("E4").Value = 0.9942998659930995
$ws.Range("F4").Value = 1.034451508610921
$ws.Range("I4").Value = 1.030088504429528
$ws.Range("J4").Value = 1.028049388238501
$ws.Range("K4").Value = 1.030915533209401
$ws.Range("L4").Value = 0.9968940712668345
$ws.Range("M4").Value = 1.03693939560252
$ws.Range("N4").Value = 1.013353593170838
# Row 5
$ws.Range("B5").Value = 1.02
$ws.Range("C5").Value = 1.023784878663995
$ws.Range("D5").Value = 1.028655637362063
$ws.Range("E5").Value = 0.9945786998346017
$ws.Range("F5").Value = 1.03482378140362
$ws.Range("I5").Value = 1.03015291564208
$ws.Range("J5").Value = 1.028277409537356
$ws.Range("K5").Value = 1.031106132170648
$ws.Range("L5").Value = 0.997117960005301
$ws.Range("M5").Value = 1.037258882212729
$ws.Range("N5").Value = 1.013430575060719
# Row 6
$ws.Range("B6").Value = 1.02
$ws.Range("C6").Value = 1.023839817004756
$ws.Range("D6").Value = 1.028696495380066
$ws.Range("E6").Value = 0.9946255319796338
$ws.Range("F6").Value = 1.034886280293858
$ws.Range("I6").Value = 1.030163706673264
$ws.Range("J6").Value = 1.028315678132047
$ws.Range("K6").Value = 1.031138114163909
$ws.Range("L6").Value = 0.9971555583673453
$ws.Range("M6").Value = 1.037312510695102
$ws.Range("N6").Value = 1.013443492627693
# Row 7
$ws.Range("B7").Value = 1.02
$ws.Range("C7").Value = 1.023461989156436
$ws.Range("D7").Value = 1.02841548038031
$ws.Range("E7").Value = 0.9943035907982488
$ws.Range("F7").Value = 1.034456483433333
$ws.Range("I7").Value = 1.030089366696817
$ws.Range("J7").Value = 1.028052436223063
$ws.Range("K7").Value = 1.030918081370697
$ws.Range("L7").Value = 0.9968970624462087
$ws.Range("M7").Value = 1.036943665587229
$ws.Range("N7").Value = 1.013354622345135
# Row 8
$ws.Range("B8").Value = 1.02
$ws.Range("C8").Value = 1.021881179039342
$ws.Range("D8").Value = 1.027239191431005
$ws.Range("E8").Value = 0.9929600610674301
$ws.Range("F8").Value = 1.032658913077265
$ws.Range("I8").Value = 1.029775175704205
$ws.Range("J8").Value = 1.026949628188569
$ws.Range("K8").Value = 1.029995416855188
$ws.Range("L8").Value = 0.995817528259106
$ws.Range("M8").Value = 1.035399799783103
$ws.Range("N8").Value = 1.012981993076818
# Row 9
$ws.Range("B9").Value = 1.02
$ws.Range("C9").Value = 1.019089541191017
$ws.Range("D9").Value = 1.025159951569914
$ws.Range("E9").Value = 0.9906006454969559
$ws.Range("F9").Value = 1.029487003785712
$ws.Range("I9").Value = 1.02920844283266
$ws.Range("J9").Value = 1.024996791095644
$ws.Range("K9").Value = 1.028358301259951
$ws.Range("L9").Value = 0.9939188001724441
$ws.Range("M9").Value = 1.032670975207532
$ws.Range("N9").Value = 1.012320940981862
# Row 10
$ws.Range("B10").Value = 1.02
$ws.Range("C10").Value = 1.017224191951641
$ws.Range("D10").Value = 1.02376935497197
$ws.Range("E10").Value = 0.989033133672735
$ws.Range("F10").Value = 1.027369225643723
$ws.Range("I10").Value = 1.02882178404928
$ws.Range("J10").Value = 1.023688359724563
$ws.Range("K10").Value = 1.02725922900441
$ws.Range("L10").Value = 0.9926553831429383
$ws.Range("M10").Value = 1.030845947054773
$ws.Range("N10").Value = 1.011877224612884
# Row 11
$ws.Range("B11").Value = 1.02
$ws.Range("C11").Value = 1.01641542215689
$ws.Range("D11").Value = 1.023166138526038
$ws.Range("E11").Value = 0.988355674866747
$ws.Range("F11").Value = 1.026451397845881
$ws.Range("I11").Value = 1.028652251840789
$ws.Range("J11").Value = 1.023120215812311
$ws.Range("K11").Value = 1.026781480629322
$ws.Range("L11").Value = 0.9921088820399291
$ws.Range("M11").Value = 1.030054264398304
$ws.Range("N11").Value = 1.011684367230116
# Row 12
$ws.Range("B12").Value = 1.02
$ws.Range("C12").Value = 1.016114845780794
$ws.Range("D12").Value = 1.022941913110113
$ws.Range("E12").Value = 0.9881042295826724
$ws.Range("F12").Value = 1.026110348399033
$ws.Range("I12").Value = 1.028588962750941
$ws.Range("J12").Value = 1.02290894134837
$ws.Range("K12").Value = 1.026603744865088
$ws.Range("L12").Value = 0.9919059725120875
$ws.Range("M12").Value = 1.029759978720007
$ws.Range("N12").Value = 1.011612621701958
# Row 13
$ws.Range("B13").Value = 1.02
$ws.Range("C13").Value = 1.01617932788341
$ws.Range("D13").Value = 1.022990017673184
$ws.Range("E13").Value = 0.9881581567098651
$ws.Range("F13").Value = 1.026183510502026
$ws.Range("I13").Value = 1.028602552846279
$ws.Range("J13").Value = 1.022954271361814
$ws.Range("K13").Value = 1.026641882434062
$ws.Range("L13").Value = 0.9919494934313052
$ws.Range("M13").Value = 1.029823113995994
$ws.Range("N13").Value = 1.011628016337917
# Row 14
$ws.Range("B14").Value = 1.02
$ws.Range("C14").Value = 1.016390579772585
$ws.Range("D14").Value = 1.023147607339027
$ws.Range("E14").Value = 0.9883348863814464
$ws.Range("F14").Value = 1.026423209217454
$ws.Range("I14").Value = 1.028647026822022
$ws.Range("J14").Value = 1.02310275672004
$ws.Range("K14").Value = 1.026766794632873
$ws.Range("L14").Value = 0.9920921077337197
$ws.Range("M14").Value = 1.030029943153859
$ws.Range("N14").Value = 1.011678438964751
# Row 15
$ws.Range("B15").Value = 1.02
$ws.Range("C15").Value = 1.016520717284999
$ws.Range("D15").Value = 1.023244681749427
$ws.Range("E15").Value = 0.9884438009545853
$ws.Range("F15").Value = 1.026570878490987
$ws.Range("I15").Value = 1.028674386639209
$ws.Range("J15").Value = 1.023194211504808
$ws.Range("K15").Value = 1.026843720168109
$ws.Range("L15").Value = 0.9921799884222134
$ws.Range("M15").Value = 1.030157348215616
$ws.Range("N15").Value = 1.011709491449334
# Row 16
$ws.Range("B16").Value = 1.02
$ws.Range("C16").Value = 1.017277844654027
$ws.Range("D16").Value = 1.023809365569488
$ws.Range("E16").Value = 0.9890781214508737
$ws.Range("F16").Value = 1.027430121256473
$ws.Range("I16").Value = 1.028832990894632
$ws.Range("J16").Value = 1.023726031971316
$ws.Range("K16").Value = 1.027290896604032
$ws.Range("L16").Value = 0.9926916645766087
$ws.Range("M16").Value = 1.030898457872393
$ws.Range("N16").Value = 1.011890008578147
# Row 17
$ws.Range("B17").Value = 1.02
$ws.Range("C17").Value = 1.017752483827563
$ws.Range("D17").Value = 1.024163286631911
$ws.Range("E17").Value = 0.9894763578477731
$ws.Range("F17").Value = 1.027968879754112
$ws.Range("I17").Value = 1.028931914635706
$ws.Range("J17").Value = 1.024059202949374
$ws.Range("K17").Value = 1.02757090386533
$ws.Range("L17").Value = 0.9930127773692701
$ws.Range("M17").Value = 1.031362949700067
$ws.Range("N17").Value = 1.012003047504502
# Row 18
$ws.Range("B18").Value = 1.02
$ws.Range("C18").Value = 1.018029230690787
$ws.Range("D18").Value = 1.024369618825822
$ws.Range("E18").Value = 0.9897087662937551
$ws.Range("F18").Value = 1.028283050073197
$ws.Range("I18").Value = 1.028989411998728
$ws.Range("J18").Value = 1.024253383231713
$ws.Range("K18").Value = 1.027734049685309
$ws.Range("L18").Value = 0.9932001317071766
$ws.Range("M18").Value = 1.031633741901466
$ws.Range("N18").Value = 1.012068911271132
# Row 19
$ws.Range("B19").Value = 1.02
$ws.Range("C19").Value = 1.018123577010103
$ws.Range("D19").Value = 1.02443995519197
$ws.Range("E19").Value = 0.9897880325774039
$ws.Range("F19").Value = 1.028390160878686
$ws.Range("I19").Value = 1.029008982674946
$ws.Range("J19").Value = 1.024319567861241
$ws.Range("K19").Value = 1.027789648112869
$ws.Range("L19").Value = 0.993264023964098
$ws.Range("M19").Value = 1.031726051757458
$ws.Range("N19").Value = 1.01209135727504
# Row 20
$ws.Range("B20").Value = 1.02
$ws.Range("C20").Value = 1.01770157012958
$ws.Range("D20").Value = 1.024125325020357
$ws.Range("E20").Value = 0.9894336180360677
$ws.Range("F20").Value = 1.027911084211017
$ws.Range("I20").Value = 1.028921322076809
$ws.Range("J20").Value = 1.024023472671082
$ws.Range("K20").Value = 1.027540880128677
$ws.Range("L20").Value = 0.9929783193494215
$ws.Range("M20").Value = 1.031313128430495
$ws.Range("N20").Value = 1.011990926736059
# Row 21
$ws.Range("B21").Value = 1.02
$ws.Range("C21").Value = 1.016328375894211
$ws.Range("D21").Value = 1.023101205651735
$ws.Range("E21").Value = 0.9882828385668249
$ws.Range("F21").Value = 1.026352627473546
$ws.Range("I21").Value = 1.028633939109696
$ws.Range("J21").Value = 1.023059038142316
$ws.Range("K21").Value = 1.026730018813613
$ws.Range("L21").Value = 0.9920501090198102
$ws.Range("M21").Value = 1.029969043214266
$ws.Range("N21").Value = 1.011663593792435
# Row 22
$ws.Range("B22").Value = 1.02
$ws.Range("C22").Value = 1.01546404767501
$ws.Range("D22").Value = 1.022456351540207
$ws.Range("E22").Value = 0.9875604150241495
$ws.Range("F22").Value = 1.025372024763645
$ws.Range("I22").Value = 1.028451413943543
$ws.Range("J22").Value = 1.022451267077506
$ws.Range("K22").Value = 1.026218584441043
$ws.Range("L22").Value = 0.9914670000341481
$ws.Range("M22").Value = 1.029122691526893
$ws.Range("N22").Value = 1.011457151461077
# Row 23
$ws.Range("B23").Value = 1.02
$ws.Range("C23").Value = 1.015922335416894
$ws.Range("D23").Value = 1.022798291548146
$ws.Range("E23").Value = 0.9879432794643023
$ws.Range("F23").Value = 1.025891932468533
$ws.Range("I23").Value = 1.028548348301949
$ws.Range("J23").Value = 1.022773590867804
$ws.Range("K23").Value = 1.026489859069259
$ws.Range("L23").Value = 0.991776070289318
$ws.Range("M23").Value = 1.02957148052458
$ws.Range("N23").Value = 1.011566650899733
# Row 24
$ws.Range("B24").Value = 1.02
$ws.Range("C24").Value = 1.017724576153583
$ws.Range("D24").Value = 1.024142478559854
$ws.Range("E24").Value = 0.9894529299347244
$ws.Range("F24").Value = 1.027937199771166
$ws.Range("I24").Value = 1.028926109026015
$ws.Range("J24").Value = 1.024039618117638
$ws.Range("K24").Value = 1.027554447111882
$ws.Range("L24").Value = 0.9929938892766442
$ws.Range("M24").Value = 1.031335640943957
$ws.Range("N24").Value = 1.01199640380564
# Row 25
$ws.Range("B25").Value = 1.02
$ws.Range("C25").Value = 1.01981198283154
$ws.Range("D25").Value = 1.025698260423169
$ws.Range("E25").Value = 0.9912096547607049
$ws.Range("F25").Value = 1.030307559771338
$ws.Range("I25").Value = 1.029356511679296
$ws.Range("J25").Value = 1.025502790146021
$ws.Range("K25").Value = 1.028782878734688
$ws.Range("L25").Value = 0.9944092447426414
$ws.Range("M25").Value = 1.033377450871438
$ws.Range("N25").Value = 1.012492367590683
